# Auto-generated Excel COM-interop script
# Applies per-cell numeric corrections to the profit/price columns (H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR leve-profit sheets,
# reflecting refreshed market-board pricing data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15:
$ws.Range("H15").Value = 4007.06
$ws.Range("I15").Value = 4007.06
$ws.Range("K15").Value = 12021.18
$ws.Range("M15").Value = -11852.18
# Row 70:
$ws.Range("H70").Value = 1910.1
$ws.Range("I70").Value = 1673.75
$ws.Range("J70").Value = 2067.6667
$ws.Range("K70").Value = 5021.25
$ws.Range("L70").Value = 6203.000100000001
$ws.Range("M70").Value = -4751.25
$ws.Range("N70").Value = -6743.000100000001
# Row 73:
$ws.Range("H73").Value = 1910.1
$ws.Range("I73").Value = 1673.75
$ws.Range("J73").Value = 2067.6667
$ws.Range("K73").Value = 5021.25
$ws.Range("L73").Value = 6203.000100000001
$ws.Range("M73").Value = -4085.25
$ws.Range("N73").Value = -8075.000100000001
# Row 86:
$ws.Range("H86").Value = 4130.2383
$ws.Range("I86").Value = 4540.4
$ws.Range("J86").Value = 4002.0625
$ws.Range("K86").Value = 4540.4
$ws.Range("L86").Value = 4002.0625
$ws.Range("M86").Value = -3417.4
$ws.Range("N86").Value = -6248.0625
# Row 89:
$ws.Range("H89").Value = 4130.2383
$ws.Range("I89").Value = 4540.4
$ws.Range("J89").Value = 4002.0625
$ws.Range("K89").Value = 22702
$ws.Range("L89").Value = 20010.3125
$ws.Range("M89").Value = -17086
$ws.Range("N89").Value = -31242.3125
# Row 100:
$ws.Range("H100").Value = 4000
$ws.Range("I100").Value = 5000
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 5000
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -4459
$ws.Range("N100").Value = -4082
# Row 107:
$ws.Range("H107").Value = 1785.9375
$ws.Range("I107").Value = 1562.7778
$ws.Range("J107").Value = 2072.8572
$ws.Range("K107").Value = 1562.7778
$ws.Range("L107").Value = 2072.8572
$ws.Range("M107").Value = 357.2221999999999
$ws.Range("N107").Value = -5912.8572
# Row 112:
$ws.Range("H112").Value = 3081.7727
$ws.Range("J112").Value = 3399.9473
$ws.Range("L112").Value = 10199.8419
$ws.Range("N112").Value = -12415.8419
# Row 116:
$ws.Range("H116").Value = 2969.973
$ws.Range("I116").Value = 2603.32
$ws.Range("K116").Value = 2603.32
$ws.Range("M116").Value = 838.6799999999998
# Row 129:
$ws.Range("H129").Value = 869.5179000000001
$ws.Range("J129").Value = 958.19147
$ws.Range("L129").Value = 2874.57441
$ws.Range("N129").Value = -12874.57441
# Row 132:
$ws.Range("H132").Value = 20844242
$ws.Range("I132").Value = 25653438
$ws.Range("K132").Value = 76960314
$ws.Range("M132").Value = -76957784
# Row 133:
$ws.Range("H133").Value = 27709
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 137:
$ws.Range("H137").Value = 1118.5938
$ws.Range("I137").Value = 996.64703
$ws.Range("J137").Value = 1256.8
$ws.Range("K137").Value = 2989.94109
$ws.Range("L137").Value = 3770.4
$ws.Range("M137").Value = -439.9410899999998
$ws.Range("N137").Value = -8870.4
# Row 138:
$ws.Range("H138").Value = 1591.2778
$ws.Range("I138").Value = 1410.4615
$ws.Range("J138").Value = 1759.1786
$ws.Range("K138").Value = 4231.3845
$ws.Range("L138").Value = 5277.5358
$ws.Range("M138").Value = 908.6154999999999
$ws.Range("N138").Value = -15557.5358

$ws = $wb.Worksheets.Item("ARM")
# Row 32:
$ws.Range("H32").Value = 3040.5857
$ws.Range("I32").Value = 2724.6516
$ws.Range("J32").Value = 8253.5
$ws.Range("K32").Value = 2724.6516
$ws.Range("L32").Value = 8253.5
$ws.Range("M32").Value = -2437.6516
$ws.Range("N32").Value = -8827.5
# Row 45:
$ws.Range("H45").Value = 1290.9445
$ws.Range("I45").Value = 1381.2142
$ws.Range("K45").Value = 1381.2142
$ws.Range("M45").Value = -1004.2142
# Row 61:
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -788
$ws.Range("N61").ClearContents()
# Row 74:
$ws.Range("H74").Value = 899.275
$ws.Range("I74").Value = 516.62964
$ws.Range("J74").Value = 1694
$ws.Range("K74").Value = 516.62964
$ws.Range("L74").Value = 1694
$ws.Range("M74").Value = 357.37036
$ws.Range("N74").Value = -3442
# Row 77:
$ws.Range("H77").Value = 899.275
$ws.Range("I77").Value = 516.62964
$ws.Range("J77").Value = 1694
$ws.Range("K77").Value = 2583.1482
$ws.Range("L77").Value = 8470
$ws.Range("M77").Value = 1784.8518
$ws.Range("N77").Value = -17206
# Row 132:
$ws.Range("H132").Value = 1808.5625
$ws.Range("I132").Value = 1485.7391
$ws.Range("K132").Value = 4457.2173
$ws.Range("M132").Value = -1927.2173
# Row 133:
$ws.Range("H133").Value = 28725.625
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 28725.625
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 28725.625
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -33785.625
# Row 136:
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -450
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 20:
$ws.Range("H20").Value = 2331
$ws.Range("I20").Value = 2218.5
$ws.Range("J20").Value = 2499.75
$ws.Range("K20").Value = 2218.5
$ws.Range("L20").Value = 2499.75
$ws.Range("M20").Value = -1971.5
$ws.Range("N20").Value = -2993.75
# Row 107:
$ws.Range("H107").Value = 1431.5883
$ws.Range("I107").Value = 934.25
$ws.Range("J107").Value = 2625.2
$ws.Range("K107").Value = 934.25
$ws.Range("L107").Value = 2625.2
$ws.Range("M107").Value = 985.75
$ws.Range("N107").Value = -6465.2
# Row 134:
$ws.Range("H134").Value = 9087.6875
$ws.Range("I134").Value = 1283.6666
$ws.Range("K134").Value = 3850.9998
$ws.Range("M134").Value = -1315.9998

$ws = $wb.Worksheets.Item("CRP")
# Row 31:
$ws.Range("H31").Value = 1684.5358
$ws.Range("I31").Value = 1265.3334
$ws.Range("J31").Value = 2168.2307
$ws.Range("K31").Value = 1265.3334
$ws.Range("L31").Value = 2168.2307
$ws.Range("M31").Value = -970.3334
$ws.Range("N31").Value = -2758.2307
# Row 34:
$ws.Range("H34").Value = 1684.5358
$ws.Range("I34").Value = 1265.3334
$ws.Range("J34").Value = 2168.2307
$ws.Range("K34").Value = 1265.3334
$ws.Range("L34").Value = 2168.2307
$ws.Range("M34").Value = -1063.3334
$ws.Range("N34").Value = -2572.2307
# Row 58:
$ws.Range("H58").Value = 2304
$ws.Range("I58").Value = 1918
$ws.Range("J58").Value = 3333.3333
$ws.Range("K58").Value = 1918
$ws.Range("L58").Value = 3333.3333
$ws.Range("M58").Value = -1715
$ws.Range("N58").Value = -3739.3333
# Row 60:
$ws.Range("H60").Value = 10129.5625
$ws.Range("I60").Value = 3235
$ws.Range("K60").Value = 3235
$ws.Range("M60").Value = -2724
# Row 132:
$ws.Range("H132").Value = 7201.5454
$ws.Range("I132").Value = 9787.538
$ws.Range("K132").Value = 29362.614
$ws.Range("M132").Value = -26832.614
# Row 134:
$ws.Range("H134").Value = 2254.3914
$ws.Range("J134").Value = 2205.6
$ws.Range("L134").Value = 6616.799999999999
$ws.Range("N134").Value = -11686.8
# Row 136:
$ws.Range("H136").Value = 2304
$ws.Range("I136").Value = 1918
$ws.Range("J136").Value = 3333.3333
$ws.Range("K136").Value = 5754
$ws.Range("L136").Value = 9999.999899999999
$ws.Range("M136").Value = -3204
$ws.Range("N136").Value = -15099.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 33:
$ws.Range("H33").Value = 205.26666
$ws.Range("I33").Value = 84.28570999999999
$ws.Range("J33").Value = 311.125
$ws.Range("K33").Value = 505.71426
$ws.Range("L33").Value = 1866.75
$ws.Range("M33").Value = -222.71426
$ws.Range("N33").Value = -2432.75
# Row 131:
$ws.Range("H131").Value = 13890080
$ws.Range("J131").Value = 1274.1212
$ws.Range("L131").Value = 3822.3636
$ws.Range("N131").Value = -13902.3636
# Row 136:
$ws.Range("H136").Value = 1492.5454
$ws.Range("I136").Value = 1390
$ws.Range("J136").Value = 1766
$ws.Range("K136").Value = 4170
$ws.Range("L136").Value = 5298
$ws.Range("M136").Value = 930
$ws.Range("N136").Value = -15498

$ws = $wb.Worksheets.Item("GSM")
# Row 132:
$ws.Range("H132").Value = 3060.5715
$ws.Range("I132").Value = 2573.6
$ws.Range("J132").Value = 4278
$ws.Range("K132").Value = 7720.799999999999
$ws.Range("L132").Value = 12834
$ws.Range("M132").Value = -5190.799999999999
$ws.Range("N132").Value = -17894

$ws = $wb.Worksheets.Item("LTW")
# Row 100:
$ws.Range("H100").Value = 2136.5
$ws.Range("I100").Value = 2038.4
$ws.Range("K100").Value = 2038.4
$ws.Range("M100").Value = -1497.4
# Row 122:
$ws.Range("H122").Value = 8069190
$ws.Range("I122").Value = 10422170
$ws.Range("J122").Value = 1830
$ws.Range("K122").Value = 31266510
$ws.Range("L122").Value = 5490
$ws.Range("M122").Value = -31264060
$ws.Range("N122").Value = -10390
# Row 132:
$ws.Range("H132").Value = 25623.928
$ws.Range("I132").Value = 1454.762
$ws.Range("J132").Value = 49793.094
$ws.Range("K132").Value = 4364.286
$ws.Range("L132").Value = 149379.282
$ws.Range("M132").Value = -1834.286
$ws.Range("N132").Value = -154439.282
# Row 136:
$ws.Range("H136").Value = 5935.15
$ws.Range("I136").Value = 8592.23
$ws.Range("K136").Value = 25776.69
$ws.Range("M136").Value = -23226.69

$ws = $wb.Worksheets.Item("WVR")
# Row 81:
$ws.Range("H81").Value = 512
$ws.Range("I81").Value = 456
$ws.Range("J81").Value = 792
$ws.Range("K81").Value = 912
$ws.Range("L81").Value = 1584
$ws.Range("M81").Value = 149
$ws.Range("N81").Value = -3706
# Row 84:
$ws.Range("H84").Value = 512
$ws.Range("I84").Value = 456
$ws.Range("J84").Value = 792
$ws.Range("K84").Value = 4560
$ws.Range("L84").Value = 7920
$ws.Range("M84").Value = 744
$ws.Range("N84").Value = -18528
# Row 96:
$ws.Range("H96").Value = 1192.1111
$ws.Range("I96").Value = 1221.1
$ws.Range("J96").Value = 1109.2858
$ws.Range("K96").Value = 1221.1
$ws.Range("L96").Value = 1109.2858
$ws.Range("M96").Value = 151.9000000000001
$ws.Range("N96").Value = -3855.2858
# Row 126:
$ws.Range("H126").Value = 52632470
$ws.Range("I126").Value = 58824340
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 176473020
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -176470550
$ws.Range("N126").Value = -9740
# Row 132:
$ws.Range("H132").Value = 2203.077
$ws.Range("I132").Value = 1849.1111
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 5547.3333
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -3017.3333
$ws.Range("N132").Value = -14058.5
# Row 136:
$ws.Range("H136").Value = 697.1429000000001
$ws.Range("I136").Value = 646.6667
$ws.Range("K136").Value = 1940.0001
$ws.Range("M136").Value = 609.9999

Write-Output "Applied Kujata_Profits market data refresh."
